# Add a new purchase record (row 15) to the shuttlecock purchase log:
# index 14, date 2024-02-19, "Bullet tournament 76", price_rod 670, n_rod 4,
# deliverly_fee 100, with the usual total / n / price_single formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 45341
$ws.Range("B15").NumberFormat = $ws.Range("B14").NumberFormat
$ws.Range("C15").Value = "Bullet tournament 76"
$ws.Range("D15").Value = 670
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 100
$ws.Range("G15").Formula = "=D15*E15+F15"
$ws.Range("H15").Formula = "=E15*12"
$ws.Range("I15").Formula = "=ROUNDUP(G15/H15,0)"

# Match the author's final selection position recorded in the sheet view.
[void]$ws.Range("C21").Select()
